{"js": "// This edit rewrites the Java stack-trace text shown in the \"idRuntimeException\"\n// sample document. The stack trace (bold, red run) moved from POI 3.16 line\n// numbers/frames to POI 3.17 ones (commit \"Fixed #253 Moving from POI 3.16 to\n// 3.17.\"). We locate the run by its distinctive first and last lines and\n// replace everything in between (inclusive) with the updated stack trace,\n// leaving the surrounding field-code runs and the following line break intact.\n\nconst body = context.document.body;\n\n// Find the very first line of the stack trace (unique in the document).\nconst startResults = body.search(\"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\", { matchCase: true, matchWholeWord: false });\nstartResults.load(\"items\");\nawait context.sync();\n\nif (startResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the stack trace start, found \" + startResults.items.length);\n}\n\n// Find the very last line of the (old) stack trace (unique in the document).\nconst endResults = body.search(\"RemoteTestRunner.main(RemoteTestRunner.java:192)\", { matchCase: true, matchWholeWord: false });\nendResults.load(\"items\");\nawait context.sync();\n\nif (endResults.items.length !== 1) {\n  throw new Error(\"Expected exactly one match for the stack trace end, found \" + endResults.items.length);\n}\n\nconst startRange = startResults.items[0].getRange(\"Start\");\nconst endRange = endResults.items[0].getRange(\"End\");\nconst fullRange = startRange.expandTo(endRange);\n\n// The updated stack trace, one array entry per line, joined with \"\\n\".\n// Note: no trailing \"\\n\" after the last line -- that character sits right\n// after fullRange's end boundary (between the last visible character and the\n// following <w:br/>), so leaving it out here avoids duplicating it.\nconst newStackTraceLines = [\n  \"divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\",\n  \"\\t/ by zero\",\n  \"java.lang.ArithmeticException: / by zero\",\n  \"\\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\",\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\",\n  \"\\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\",\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\",\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\",\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\",\n  \"\\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\",\n  \"\\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1074)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1)\",\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:160)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\",\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\",\n  \"\\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\",\n  \"\\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\",\n  \"\\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\",\n  \"\\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\",\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\",\n  \"\\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\",\n  \"\\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\",\n  \"\\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\",\n  \"\\tat java.lang.reflect.Method.invoke(Method.java:498)\",\n  \"\\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\",\n  \"\\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\",\n  \"\\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\",\n  \"\\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\",\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\",\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\",\n  \"\\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\",\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\",\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\",\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\",\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\",\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\",\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\",\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\",\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\",\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\",\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\",\n  \"\\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:128)\",\n  \"\\tat org.junit.runners.Suite.runChild(Suite.java:27)\",\n  \"\\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\",\n  \"\\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\",\n  \"\\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\",\n  \"\\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\",\n  \"\\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\",\n  \"\\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\",\n  \"\\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\",\n  \"\\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\"\n];\nconst newStackTrace = newStackTraceLines.join(\"\\n\");\n\nfullRange.insertText(newStackTrace, \"Replace\");\n\nawait context.sync();\n", "ps1": "# This edit rewrites the Java stack-trace text shown in the \"idRuntimeException\"\n# sample document. The stack trace (bold, red run) moved from POI 3.16 line\n# numbers/frames to POI 3.17 ones (commit \"Fixed #253 Moving from POI 3.16 to\n# 3.17.\"). We locate the run by its distinctive first and last lines and\n# replace everything in between (inclusive) with the updated stack trace,\n# leaving the surrounding field-code runs and the following line break intact.\n\n$d = $word.ActiveDocument\n\n# Locate the first line of the stack trace.\n$startRng = $d.Content\n$startRng.Find.ClearFormatting()\n$foundStart = $startRng.Find.Execute('divOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:')\nif (-not $foundStart) {\n    throw \"Could not find the start of the stack trace.\"\n}\n\n# Locate the last line of the (old) stack trace.\n$endRng = $d.Content\n$endRng.Find.ClearFormatting()\n$foundEnd = $endRng.Find.Execute('RemoteTestRunner.main(RemoteTestRunner.java:192)')\nif (-not $foundEnd) {\n    throw \"Could not find the end of the stack trace.\"\n}\n\n# Build a range spanning the whole stack trace (start of first line through\n# end of last line) and overwrite its text in one shot, preserving the\n# run's existing bold/red formatting as well as the surrounding fldChar\n# field-code runs and the trailing line break.\n$full = $d.Range($startRng.Start, $endRng.End)\n\n$newStackTrace = @'\ndivOp(java.lang.Integer,java.lang.Integer) with arguments [1, 0] failed:\n\t/ by zero\njava.lang.ArithmeticException: / by zero\n\tat org.eclipse.acceleo.query.services.NumberServices.divOp(NumberServices.java:99)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke0(Native Method)\n\tat sun.reflect.NativeMethodAccessorImpl.invoke(NativeMethodAccessorImpl.java:62)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.eclipse.acceleo.query.runtime.impl.JavaMethodService.internalInvoke(JavaMethodService.java:163)\n\tat org.eclipse.acceleo.query.runtime.impl.AbstractService.invoke(AbstractService.java:136)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.callService(EvaluationServices.java:129)\n\tat org.eclipse.acceleo.query.runtime.impl.EvaluationServices.call(EvaluationServices.java:168)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.caseCall(AstEvaluator.java:189)\n\tat org.eclipse.acceleo.query.ast.util.AstSwitch.doSwitch(AstSwitch.java:118)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.eclipse.acceleo.query.parser.AstEvaluator.eval(AstEvaluator.java:112)\n\tat org.eclipse.acceleo.query.runtime.impl.QueryEvaluationEngine.eval(QueryEvaluationEngine.java:52)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1074)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseUserDoc(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:160)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1254)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseBlock(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:199)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:275)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.caseDocumentTemplate(M2DocEvaluator.java:1)\n\tat org.obeonetwork.m2doc.template.util.TemplateSwitch.doSwitch(TemplateSwitch.java:279)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:53)\n\tat org.eclipse.emf.ecore.util.Switch.doSwitch(Switch.java:69)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.doSwitch(M2DocEvaluator.java:1038)\n\tat org.obeonetwork.m2doc.generator.M2DocEvaluator.generate(M2DocEvaluator.java:264)\n\tat org.obeonetwork.m2doc.util.M2DocUtils.generate(M2DocUtils.java:712)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.prepareoutputAndGenerate(AbstractTemplatesTestSuite.java:459)\n\tat org.obeonetwork.m2doc.tests.AbstractTemplatesTestSuite.generation(AbstractTemplatesTestSuite.java:369)\n\tat sun.reflect.GeneratedMethodAccessor75.invoke(Unknown Source)\n\tat sun.reflect.DelegatingMethodAccessorImpl.invoke(DelegatingMethodAccessorImpl.java:43)\n\tat java.lang.reflect.Method.invoke(Method.java:498)\n\tat org.junit.runners.model.FrameworkMethod$1.runReflectiveCall(FrameworkMethod.java:50)\n\tat org.junit.internal.runners.model.ReflectiveCallable.run(ReflectiveCallable.java:12)\n\tat org.junit.runners.model.FrameworkMethod.invokeExplosively(FrameworkMethod.java:47)\n\tat org.junit.internal.runners.statements.InvokeMethod.evaluate(InvokeMethod.java:17)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.runLeaf(ParentRunner.java:325)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:78)\n\tat org.junit.runners.BlockJUnit4ClassRunner.runChild(BlockJUnit4ClassRunner.java:57)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.internal.runners.statements.RunAfters.evaluate(RunAfters.java:27)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.junit.runners.Suite.runChild(Suite.java:128)\n\tat org.junit.runners.Suite.runChild(Suite.java:27)\n\tat org.junit.runners.ParentRunner$3.run(ParentRunner.java:290)\n\tat org.junit.runners.ParentRunner$1.schedule(ParentRunner.java:71)\n\tat org.junit.runners.ParentRunner.runChildren(ParentRunner.java:288)\n\tat org.junit.runners.ParentRunner.access$000(ParentRunner.java:58)\n\tat org.junit.runners.ParentRunner$2.evaluate(ParentRunner.java:268)\n\tat org.junit.runners.ParentRunner.run(ParentRunner.java:363)\n\tat org.eclipse.jdt.internal.junit4.runner.JUnit4TestReference.run(JUnit4TestReference.java:86)\n\tat org.eclipse.jdt.internal.junit.runner.TestExecution.run(TestExecution.java:38)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:539)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.runTests(RemoteTestRunner.java:761)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.run(RemoteTestRunner.java:461)\n\tat org.eclipse.jdt.internal.junit.runner.RemoteTestRunner.main(RemoteTestRunner.java:207)\n'@\n\n$full.Text = $newStackTrace\n"}
